$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to reflect that errors/invalids are now percentages
$ws.Range("AZ1").Value = "accuracy_SEVERITY_ERRORS"
$ws.Range("BA1").Value = "accuracy_CWE_ERRORS"
$ws.Range("BB1").Value = "accuracy_INVALID_CWE_INFERENCE_counter"
$ws.Range("BC1").Value = "accuracy_INVALID_SEVERITY_INFERENCE_counter"

# Convert raw error counts (AZ/BA columns) into percentages by dividing by 5
$lastRow = $ws.Cells.Item($ws.Rows.Count, 52).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $azCell = $ws.Cells.Item($r, 52)  # column AZ
    $baCell = $ws.Cells.Item($r, 53)  # column BA

    $azVal = $azCell.Value2
    $baVal = $baCell.Value2

    if ($azVal -ne $null) {
        $azCell.Value = $azVal / 5
    }
    if ($baVal -ne $null) {
        $baCell.Value = $baVal / 5
    }
}
